$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.160.56"
$ws.Range("E2").Value = "  -1.65%  "

$ws.Range("D3").Value = "1.797.33"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.99"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.552"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.24"
$ws.Range("E8").Value = "  -0.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.286"
$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0721"
$ws.Range("E10").Value = "  +5.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("E11").Value = "  -1.31%  "

$ws.Range("D12").Value = "2.055.00"

$ws.Range("D13").Value = "1.802.93"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.74"
$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "34.123.17"
$ws.Range("E16").Value = "  -1.68%  "

$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.17"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.79"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  +3.28%  "

$ws.Range("E23").Value = "  -1.74%  "

$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.93"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.54"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0522"
$ws.Range("E30").Value = "  +1.28%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.72"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").Value = "1.414.60"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("E36").Value = "  +2.36%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.946"
$ws.Range("E39").Value = "  +4.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.29"
$ws.Range("E40").Value = "  -3.03%  "

$ws.Range("E41").Value = "  -2.60%  "

$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("E43").Value = "  +4.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.96"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.74"
$ws.Range("E46").Value = "  +2.26%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.953.54"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.04"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.92"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("E51").Value = "  -0.17%  "

